$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.463.30'
$ws.Cells.Item(2, 5).Value = '  +0.16%  '
$ws.Cells.Item(3, 4).Value = '1.865.84'
$ws.Cells.Item(3, 5).Value = '  -0.46%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).Value = '''235.25'
$ws.Cells.Item(5, 5).Value = '  -1.11%  '
$ws.Cells.Item(6, 4).Value = '''1.000'
$ws.Cells.Item(6, 5).Value = '  +0.06%  '
$ws.Cells.Item(7, 4).Value = '''0.4829'
$ws.Cells.Item(7, 5).Value = '  +0.17%  '
$ws.Cells.Item(8, 4).Value = '''0.2805'
$ws.Cells.Item(8, 5).Value = '  -0.75%  '
$ws.Cells.Item(9, 4).Value = '''0.06505'
$ws.Cells.Item(9, 5).Value = '  -0.78%  '
$ws.Cells.Item(10, 4).Value = '1.881.49'
$ws.Cells.Item(10, 5).Value = '  +0.48%  '
$ws.Cells.Item(11, 4).Value = '''0.07420'
$ws.Cells.Item(11, 5).Value = '  -0.53%  '
$ws.Cells.Item(12, 4).Value = '''16.36'
$ws.Cells.Item(12, 5).Value = '  -0.33%  '
$ws.Cells.Item(13, 4).Value = '''5.062'
$ws.Cells.Item(13, 5).Value = '  -0.55%  '
$ws.Cells.Item(14, 4).Value = '''87.20'
$ws.Cells.Item(14, 5).Value = '  -1.19%  '
$ws.Cells.Item(15, 4).Value = '''0.6468'
$ws.Cells.Item(16, 4).Value = '30.445.03'
$ws.Cells.Item(16, 5).Value = '  +0.36%  '
$ws.Cells.Item(17, 4).Value = '''1.000'
$ws.Cells.Item(17, 5).Value = '  +0.04%  '
$ws.Cells.Item(18, 4).Value = '''12.98'
$ws.Cells.Item(18, 5).Value = '  -2.51%  '
$ws.Cells.Item(19, 4).Value = '''234.14'
$ws.Cells.Item(19, 5).Value = '  +5.71%  '
$ws.Cells.Item(20, 4).Value = '''0.000007540'
$ws.Cells.Item(20, 5).Value = '  -1.42%  '
$ws.Cells.Item(21, 4).Value = '2.118.54'
$ws.Cells.Item(21, 5).Value = '  -0.05%  '
$ws.Cells.Item(22, 5).Value = '  +0.08%  '
$ws.Cells.Item(23, 4).Value = '''5.149'
$ws.Cells.Item(23, 5).Value = '  -3.29%  '
$ws.Cells.Item(24, 4).Value = '''6.090'
$ws.Cells.Item(24, 5).Value = '  -1.54%  '
$ws.Cells.Item(25, 4).Value = '''9.327'
$ws.Cells.Item(25, 5).Value = '  +0.68%  '
$ws.Cells.Item(26, 4).Value = '''166.95'
$ws.Cells.Item(26, 5).Value = '  +0.88%  '
$ws.Cells.Item(27, 4).Value = '''18.36'
$ws.Cells.Item(27, 5).Value = '  -1.40%  '
$ws.Cells.Item(28, 4).Value = '''1.923'
$ws.Cells.Item(28, 5).Value = '  -2.05%  '
$ws.Cells.Item(29, 4).Value = '''0.1024'
$ws.Cells.Item(29, 5).Value = '  +9.07%  '
$ws.Cells.Item(30, 4).Value = '''1.375'
$ws.Cells.Item(30, 5).Value = '  -5.49%  '
$ws.Cells.Item(31, 4).Value = '''4.272'
$ws.Cells.Item(31, 5).Value = '  -0.88%  '
$ws.Cells.Item(32, 4).Value = '''3.999'
$ws.Cells.Item(32, 5).Value = '  -0.62%  '
$ws.Cells.Item(33, 4).Value = '''0.04984'
$ws.Cells.Item(33, 5).Value = '  -1.41%  '
$ws.Cells.Item(34, 4).Value = '''1.177'
$ws.Cells.Item(34, 5).Value = '  -2.99%  '
$ws.Cells.Item(35, 4).Value = '''0.7313'
$ws.Cells.Item(35, 5).Value = '  -3.37%  '
$ws.Cells.Item(36, 4).Value = '''0.9998'
$ws.Cells.Item(36, 5).Value = '  +0.19%  '
$ws.Cells.Item(37, 4).Value = '''2.711'
$ws.Cells.Item(37, 5).Value = '  -0.13%  '
$ws.Cells.Item(38, 4).Value = '''0.01917'
$ws.Cells.Item(38, 5).Value = '  +4.37%  '
$ws.Cells.Item(39, 4).Value = '''2.633'
$ws.Cells.Item(39, 5).Value = '  +0.35%  '
$ws.Cells.Item(40, 4).Value = '''0.9131'
$ws.Cells.Item(40, 5).Value = '  +0.85%  '
$ws.Cells.Item(41, 4).Value = '''2.046'
$ws.Cells.Item(41, 5).Value = '  -1.79%  '
$ws.Cells.Item(42, 4).Value = '''106.17'
$ws.Cells.Item(42, 5).Value = '  -0.57%  '
$ws.Cells.Item(43, 4).Value = '''0.9957'
$ws.Cells.Item(43, 5).Value = '  -0.73%  '
$ws.Cells.Item(44, 4).Value = '''0.4204'
$ws.Cells.Item(44, 5).Value = '  -2.12%  '
$ws.Cells.Item(45, 4).Value = '''5.563'
$ws.Cells.Item(45, 5).Value = '  -6.35%  '
$ws.Cells.Item(46, 4).Value = '''7.230'
$ws.Cells.Item(46, 5).Value = '  -3.24%  '
$ws.Cells.Item(47, 4).Value = '''61.99'
$ws.Cells.Item(47, 5).Value = '  -6.55%  '
$ws.Cells.Item(48, 4).Value = '''0.1228'
$ws.Cells.Item(48, 5).Value = '  -5.59%  '
$ws.Cells.Item(49, 4).Value = '''8.851'
$ws.Cells.Item(49, 5).Value = '  -0.92%  '
$ws.Cells.Item(50, 4).Value = '''1.441'
$ws.Cells.Item(50, 5).Value = '  -2.79%  '
$ws.Cells.Item(51, 4).Value = '''33.60'
$ws.Cells.Item(51, 5).Value = '  -1.90%  '
